# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (column G) values calculated from s_vals regen, replacing the old Strike# values
$kValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 2
    7  = 0
    8  = 0
    9  = 1
    10 = 1
    11 = 2
    12 = 2
    13 = 0
    14 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
